$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet Data")

# Add new worksheet "ABC" right after "Sheet Data"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "ABC"

# Header row
$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Age"
$ws2.Range("C1").Value = "Email"
$ws2.Range("D1").Value = "Mobile"
$ws2.Range("E1").Value = "Gender"
$ws2.Range("F1").Value = "Salary"
$ws2.Range("G1").Value = "Address"

# Data row
$ws2.Range("A2").Value = "Rohit"
$ws2.Range("B2").Value = 26
$ws2.Range("C2").Value = "rohit@mail"
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:rohit@mail", "", "", "rohit@mail")
$ws2.Range("C2").Style = "Hyperlink"
$ws2.Range("D2").Value = 398539
$ws2.Range("E2").Value = "M"
$ws2.Range("F2").Value = 3000
$ws2.Range("G2").Value = "Noida"

# Set selections to match final state
$ws1.Range("C8").Select() | Out-Null
$ws2.Range("G5").Select() | Out-Null
$ws2.Activate() | Out-Null
